$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.791792750358582
$ws.Range("B1").Value = 4.292407989501953
$ws.Range("C1").Value = 1.587259411811829
$ws.Range("D1").Value = 0.8649300932884216
$ws.Range("E1").Value = 0.4668846726417542
